# Apply the timetable update: adjust hour allotments for rows 3-5,
# replace the "IV CSE" rows (8-10) and the existing "VI CSE" rows (11-12)
# with the finalized VI CSE curriculum, and append the newly added
# VI CSE rows (13-17: value-added course, spoken tutorial, the two
# integrated-lab rows, and the mini project).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some cells need to become a *blank text* value (stored as an
# empty string, matching the original sheet's blank "str" cells) rather
# than a fully-cleared/number cell. Assigning a bare apostrophe forces
# Excel's text interpretation (yielding an empty string), then resetting
# the cell style back to Normal drops the quote-prefix formatting that
# the apostrophe entry would otherwise leave behind.
function Set-BlankText($rng) {
    $rng.Value = "'"
    $rng.Style = "Normal"
}

# --- Row 3: Probability, Statistics and Linear Algebra — hours allotted 5 -> 3
$ws.Range("L3").Value = 3

# --- Row 4: Quantum Physics — hours allotted 4 -> 3
$ws.Range("L4").Value = 3

# --- Row 5: Programming in C — hours allotted 4 -> 3
$ws.Range("L5").Value = 3

# --- Row 8: IV CSE / Theory of Computation -> VI CSE / Cryptography and cyber security (Integrated)
$ws.Range("A8").Value = "VI CSE"
$ws.Range("B8").Value = "CS2611"
$ws.Range("C8").Value = "Cryptography and cyber security (Integrated)"
$ws.Range("D8").Value = "MST"
$ws.Range("E8").Value = "GM"
$ws.Range("F8").Value = "VP"
$ws.Range("G8").Value = "ND"
$ws.Range("M8").Value = 0

# --- Row 9: IV CSE / Database Management Systems -> VI CSE / Internet of Things (Integrated)
$ws.Range("A9").Value = "VI CSE"
$ws.Range("B9").Value = "CS2612"
$ws.Range("C9").Value = "Internet of Things (Integrated)"
$ws.Range("D9").Value = "ATP"
$ws.Range("E9").Value = "SSB"
$ws.Range("F9").Value = "CS"
$ws.Range("G9").Value = "MJ"
$ws.Range("L9").Value = 3

# --- Row 10: IV CSE / Computer Networks -> VI CSE / Software Defined Networks (Open Elective I)
$ws.Range("A10").Value = "VI CSE"
$ws.Range("B10").Value = "EC2014"
$ws.Range("C10").Value = "Software Defined Networks -Open Elective - I*"
$ws.Range("D10").Value = "SD1"
Set-BlankText($ws.Range("E10"))
$ws.Range("F10").Value = "SD2"
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = "ECE"
$ws.Range("L10").Value = 3

# --- Row 11: Cryptography and Cyber Security -> Renewable Energy Technologies (Open Elective I)
$ws.Range("B11").Value = "ME2011"
$ws.Range("C11").Value = "Renewable Energy Technologies -Open Elective - I*"
$ws.Range("D11").Value = "RET1"
Set-BlankText($ws.Range("E11"))
$ws.Range("F11").Value = "RET2"
Set-BlankText($ws.Range("G11"))
$ws.Range("I11").Value = 2
$ws.Range("J11").Value = "MECH"
$ws.Range("L11").Value = 3
$ws.Range("M11").Value = 1

# --- Row 12: Internet of Things -> Image Processing (Professional Elective III)
$ws.Range("B12").Value = "CS2V62"
$ws.Range("C12").Value = "Image Processing - Professional Elective - III"
$ws.Range("D12").Value = "DMDP"
$ws.Range("E12").Value = "VNK"
Set-BlankText($ws.Range("F12"))
Set-BlankText($ws.Range("G12"))
$ws.Range("I12").Value = 2
$ws.Range("L12").Value = 4

# --- Row 13 (new): Value added course - Entrepreneurship Development
$ws.Range("A13").Value = "VI CSE"
$ws.Range("B13").Value = "ED2VA1"
$ws.Range("C13").Value = "Value added course - Entrepreneurship Development"
$ws.Range("D13").Value = "PAC"
$ws.Range("E13").Value = "NPP"
$ws.Range("F13").Value = "SGR"
$ws.Range("G13").Value = "RSA"
Set-BlankText($ws.Range("H13"))
$ws.Range("I13").Value = 4
$ws.Range("J13").Value = "CSE"
$ws.Range("K13").Value = "NO"
$ws.Range("L13").Value = 4
$ws.Range("M13").Value = 0

# --- Row 14 (new): IIT Spoken Tutorial class
$ws.Range("A14").Value = "VI CSE"
$ws.Range("B14").Value = "IT3412"
$ws.Range("C14").Value = "IIT Spoken Tutorial class"
$ws.Range("D14").Value = "LA1"
$ws.Range("E14").Value = "LA2"
$ws.Range("F14").Value = "LA3"
$ws.Range("G14").Value = "LA4"
Set-BlankText($ws.Range("H14"))
$ws.Range("I14").Value = 4
$ws.Range("J14").Value = "IT"
$ws.Range("K14").Value = "NO"
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 1

# --- Row 15 (new): Cryptography and cyber security (Integrated Lab)
$ws.Range("A15").Value = "VI CSE"
$ws.Range("B15").Value = "CS2611"
$ws.Range("C15").Value = "Cryptography and cyber security (Integrated Lab)"
$ws.Range("D15").Value = "MST"
$ws.Range("E15").Value = "GM"
$ws.Range("F15").Value = "VP"
$ws.Range("G15").Value = "ND"
Set-BlankText($ws.Range("H15"))
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = "CSE"
$ws.Range("K15").Value = "NO"
$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 0

# --- Row 16 (new): Internet of Things (Integrated Lab)
$ws.Range("A16").Value = "VI CSE"
$ws.Range("B16").Value = "CS2612"
$ws.Range("C16").Value = "Internet of Things (Integrated Lab)"
$ws.Range("D16").Value = "ATP"
$ws.Range("E16").Value = "SSB"
$ws.Range("F16").Value = "CS"
$ws.Range("G16").Value = "MJ"
Set-BlankText($ws.Range("H16"))
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = "CSE"
$ws.Range("K16").Value = "NO"
$ws.Range("L16").Value = 2
$ws.Range("M16").Value = 0

# --- Row 17 (new): MiniProject
$ws.Range("A17").Value = "VI CSE"
$ws.Range("B17").Value = "CS2698"
$ws.Range("C17").Value = "MiniProject"
$ws.Range("D17").Value = "SAA"
$ws.Range("E17").Value = "RAS"
$ws.Range("F17").Value = "RSK"
$ws.Range("G17").Value = "RSA"
Set-BlankText($ws.Range("H17"))
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = "CSE"
$ws.Range("K17").Value = "NO"
$ws.Range("L17").Value = 3
$ws.Range("M17").Value = 0
